$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "parent_directory"
$ws.Range("B21").Value = "~/Programming/github/HydroTroPe"

$ws.Range("A22").Value = "data_parent_folder"
$ws.Range("B22").Value = "~/Dropbox/PhD/Computation/ForestCarbon/2022 Kalimantan customer work/0. Raw Data"

$ws.Range("A23").Value = "output_directory"
$ws.Range("B23").Value = "output"
$ws.Range("C23").Value = "relative to parent_directory"

$ws.Range("A24").Select()
